$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95 (shifts existing rows 95..182 down to 96..183)
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new data record
$ws.Range("A95").Value2 = 6
$ws.Range("B95").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C95").Value2 = "Metropolitana"
$ws.Range("D95").Value2 = 44494
$ws.Range("E95").Value2 = 13
$ws.Range("F95").Value2 = "Fruta"
$ws.Range("G95").Value2 = 100101
$ws.Range("H95").Value2 = "Berries"
$ws.Range("I95").Value2 = 100101001
$ws.Range("J95").Value2 = "Arándano (blue)"
$ws.Range("K95").Value2 = "Sin especificar"
$ws.Range("L95").Value2 = "Especial"
$ws.Range("M95").Value2 = 150
$ws.Range("N95").Value2 = 12000
$ws.Range("O95").Value2 = 12000
$ws.Range("P95").Value2 = 12000
$ws.Range("Q95").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R95").Value2 = "Provincia del Elquí"
$ws.Range("S95").Value2 = 6000
$ws.Range("T95").Value2 = 2
